$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.629
$ws.Range("C2").Value = 0.10422342

$ws.Range("B3").Value = 1.32366666666667
$ws.Range("C3").Value = 0.09395386

$ws.Range("B4").Value = 1.45466666666667
$ws.Range("C4").Value = 0.0985100266666667

$ws.Range("B5").Value = 4.54333333333333
$ws.Range("C5").Value = 0.174100533333333

$ws.Range("B6").Value = 4.65733333333333
$ws.Range("C6").Value = 0.17632664

$ws.Range("B7").Value = 4.833
$ws.Range("C7").Value = 0.17959428

$ws.Range("B8").Value = 18.8603333333333
$ws.Range("C8").Value = 0.35419706

$ws.Range("B9").Value = 18.986
$ws.Range("C9").Value = 0.35579764

$ws.Range("B10").Value = 19.246
$ws.Range("C10").Value = 0.3579756

$ws.Range("B11").Value = 37.1543333333333
$ws.Range("C11").Value = 0.497868066666667

$ws.Range("B12").Value = 37.1336666666667
$ws.Range("C12").Value = 0.498333806666667

$ws.Range("B13").Value = 36.2533333333333
$ws.Range("C13").Value = 0.4915952

$ws.Range("B14").Value = 92.3036666666667
$ws.Range("C14").Value = 0.784581166666667

$ws.Range("B15").Value = 89.988
$ws.Range("C15").Value = 0.7738968

$ws.Range("B16").Value = 92.5273333333333
$ws.Range("C16").Value = 0.784631786666667

$ws.Range("B17").Value = 176.261666666667
$ws.Range("C17").Value = 1.08224663333333

$ws.Range("B18").Value = 175.535666666667
$ws.Range("C18").Value = 1.08481042

$ws.Range("B19").Value = 173.762
$ws.Range("C19").Value = 1.0773244
